$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (regen sval data to filter save games)
$ws.Range("B2").Value = 0.003208871385164791
$ws.Range("C2").Value = 0.0000005461030343489881
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1133.79286921213

# Update row 3 values
$ws.Range("B3").Value = 0.2917716402565462
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 1135.737209517158
